$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - symbol changed to MATIC
$ws.Range("A2").Value = "BINANCE_SPOT_MATIC_USDT"
$ws.Range("C2").Value = "0/0    ↗"
$ws.Range("D2").Value = "2/0    ↗"
$ws.Range("E2").Value = "0/0    ↘"

# Row 3 - symbol changed to LUNA
$ws.Range("A3").Value = "BINANCE_SPOT_LUNA_USDT"
$ws.Range("B3").Value = "0/0    ↘"
$ws.Range("C3").Value = "0/0    ↗"
$ws.Range("D3").Value = "0/0    ↗"
$ws.Range("E3").Value = "0/0    ↗"
$ws.Range("F3").Value = "1/0    ↘"

# Row 4 - ETH (symbol unchanged)
$ws.Range("B4").Value = "0/0    ↗"
$ws.Range("C4").Value = "0/0    ↗"
$ws.Range("D4").Value = "2/0    ↘"
$ws.Range("E4").Value = "2/0    ↘"

# Row 5 - ADA (symbol unchanged)
$ws.Range("B5").Value = "2/0    ↗"
$ws.Range("C5").Value = "0/0    ↗"
$ws.Range("E5").Value = "0/0    ↘"
$ws.Range("F5").Value = "1/0    ↘"

# Row 6 - SOL (symbol unchanged)
$ws.Range("B6").Value = "1/0    ↘"
$ws.Range("D6").Value = "0/0    ↘"
$ws.Range("E6").Value = "0/0    ↘"
$ws.Range("F6").Value = "0/0    ↘"

# Row 7 - CRV (symbol unchanged)
$ws.Range("E7").Value = "1/0    ↘"
$ws.Range("F7").Value = "1/0    ↘"

# Row 8 - BTC (symbol unchanged)
$ws.Range("B8").Value = "4/1    ↘"
$ws.Range("C8").Value = "0/0    ↗"
$ws.Range("D8").Value = "0/0    ↘"
$ws.Range("E8").Value = "1/0    ↘"
$ws.Range("F8").Value = "1/0    ↘"
